# Update the cached "datetimeFigureOut" date field text from 12-01-2020 to
# 14-01-2020 everywhere it is stored: the Slide Master and every Slide
# Layout (CustomLayout) date placeholder.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "12-01-2020") {
                $sh.TextFrame.TextRange.Text = "14-01-2020"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout.Shapes
}
